$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.329.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.578.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.47"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.802.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.582.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.517"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "28.387.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "231.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0687"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.104"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0481"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.55%  "
$ws.Range("E32").Value = "  -2.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.392.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  +7.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.71%  "
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("E39").Value = "  +3.73%  "
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.519"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.786"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.928"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "62.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.715.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "85.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "41.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.32%  "
